# Reverse the order of the comma-separated "Recorded By" entries in column G
# for every data row of the active sheet (row 1 is the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2

    if ($current -ne $null -and $current -ne "") {
        $parts = $current -split ", "
        if ($parts.Length -gt 1) {
            $reversed = $parts[($parts.Length - 1)..0]
            $cell.Value = [string]::Join(", ", $reversed)
        }
    }
}
